# Update "想去人数" (Number of people interested) values in column F
# for the "展览" (Exhibition) and "全部类型" (All Types) worksheets,
# matching the refreshed data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2243
$ws.Range("F3").Value = 123
$ws.Range("F4").Value = 74
$ws.Range("F5").Value = 705
$ws.Range("F9").Value = 2662
$ws.Range("F10").Value = 1634
$ws.Range("F11").Value = 1666
$ws.Range("F13").Value = 274
$ws.Range("F14").Value = 689
$ws.Range("F15").Value = 851
$ws.Range("F16").Value = 128
$ws.Range("F17").Value = 350
$ws.Range("F18").Value = 1104
$ws.Range("F21").Value = 532
$ws.Range("F22").Value = 5906
$ws.Range("F24").Value = 1137
$ws.Range("F26").Value = 169
$ws.Range("F28").Value = 273
$ws.Range("F32").Value = 852
$ws.Range("F34").Value = 75
$ws.Range("F35").Value = 241
$ws.Range("F36").Value = 434
$ws.Range("F37").Value = 1241
$ws.Range("F41").Value = 134
$ws.Range("F42").Value = 142

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2243
$ws.Range("F5").Value = 123
$ws.Range("F6").Value = 74
$ws.Range("F7").Value = 705
$ws.Range("F14").Value = 2662
$ws.Range("F15").Value = 1634
$ws.Range("F16").Value = 1666
$ws.Range("F18").Value = 274
$ws.Range("F19").Value = 689
$ws.Range("F21").Value = 851
$ws.Range("F22").Value = 128
$ws.Range("F23").Value = 350
$ws.Range("F24").Value = 1104
$ws.Range("F26").Value = 532
$ws.Range("F27").Value = 5906
$ws.Range("F29").Value = 1137
$ws.Range("F31").Value = 169
$ws.Range("F33").Value = 273
$ws.Range("F37").Value = 852
$ws.Range("F39").Value = 75
$ws.Range("F40").Value = 434
$ws.Range("F41").Value = 1241
$ws.Range("F45").Value = 134
$ws.Range("F46").Value = 142
